# Updated model comparison to include survreg
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Survreg" column (column E) values
$ws.Range("E1").Value = "Survreg"
$ws.Range("E1").Font.Bold = $true

$ws.Range("E2").Value = "6.0 (4.7, 8.0)"
$ws.Range("E3").Value = 1.35
$ws.Range("E5").Value = 3.7
$ws.Range("E6").Value = 10.1

$ws.Range("E8").Value = "12.15 (7.9, 18.6)"
$ws.Range("E9").Value = "2.5 (passed in to get it to fit)"
$ws.Range("E10").Value = 2.7
$ws.Range("E11").Value = 10.5
$ws.Range("E12").Value = 18.8

$ws.Range("E14").Value = "10.3 (7.0,15.4)"
$ws.Range("E15").Value = "3.8 (passed in to get fit)"
$ws.Range("E16").Value = 4.7
$ws.Range("E17").Value = 9.4
$ws.Range("E18").Value = 13.8

# Resize the new column to fit its contents
$ws.Columns.Item(5).AutoFit()

# Update the saved selection/cursor position
[void]$ws.Range("E19").Select()

# Match the page orientation change recorded for the sheet
$ws.PageSetup.Orientation = 1
